# Doing Updates for Financials
#
# The DAR sheet grows a new first data-column (new column D, with period
# ending 2019-01-31 / serial 43463). All former columns D:K (each period's
# Income Statement / Balance Sheet / Cash Flow figures) shift right to E:L,
# and a handful of the shifted figures are themselves restated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column before D; this shifts the existing D:K data to E:L
# (values + styles move together).
$ws.Columns("D").Insert()

# The freshly inserted column D inherits column C's style. Re-stripe it from
# column E (its new right-hand neighbor, i.e. the old column D) so the date
# row keeps the date format and the numeric rows keep the number format.
$ws.Range("E5:E102").Copy()
$ws.Range("D5:D102").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---- Populate the new column D (the new "most recent period") ----
$ws.Range("D7").Value = 43463
$ws.Range("D8").Value = 3387700
$ws.Range("D9").Value = 2649500
$ws.Range("D10").Value = 738200
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 51000
$ws.Range("D15").Value = 321200
$ws.Range("D17").Value = 3328600
$ws.Range("D18").Value = 59200
$ws.Range("D20").Value = 145200
$ws.Range("D21").Value = 525600
$ws.Range("D22").Value = 86400
$ws.Range("D23").Value = 118000
$ws.Range("D24").Value = 13700
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 104200
$ws.Range("D27").Value = 99800
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = 1700
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = -145200
$ws.Range("D33").Value = 101500
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 101500

$ws.Range("D38").Value = 43463
$ws.Range("D41").Value = 107300
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 392200
$ws.Range("D44").Value = 341000
$ws.Range("D45").Value = 57500
$ws.Range("D46").Value = 897900
$ws.Range("D47").Value = 410200
$ws.Range("D48").Value = 1687900
$ws.Range("D49").Value = 1825000
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 68400
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 4889400
$ws.Range("D57").Value = 219500
$ws.Range("D58").Value = 7500
$ws.Range("D59").Value = 313500
$ws.Range("D60").Value = 540500
$ws.Range("D61").Value = 1666900
$ws.Range("D62").Value = 346100
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 2616300
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 1087500
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 2273000
$ws.Range("D77").Value = 0

$ws.Range("D80").Value = 43463
$ws.Range("D81").Value = 101500
$ws.Range("D83").Value = 321200
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 398600
$ws.Range("D91").Value = -321900
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -342400
$ws.Range("D96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = -47600
$ws.Range("D101").Value = -8200
$ws.Range("D102").Value = 500

# ---- A handful of the shifted-over figures (now in E:K, formerly D:J) were
# also restated and need correcting beyond the plain shift ----
$ws.Range("E9").Value = 5758200
$ws.Range("F9").Value = 2635300
$ws.Range("E10").Value = -2095900
$ws.Range("F10").Value = 756600
$ws.Range("E14").Value = 900
$ws.Range("E17").Value = 3521900
$ws.Range("F17").Value = 3237200
$ws.Range("E18").Value = 140300
$ws.Range("F18").Value = 154700
$ws.Range("E20").Value = 12800
$ws.Range("F20").Value = 62000
$ws.Range("E32").Value = -12800
$ws.Range("F32").Value = -62000
$ws.Range("G59").Value = 486300
$ws.Range("E89").Value = 410400
$ws.Range("E91").Value = -274200
$ws.Range("G91").Value = -229800
$ws.Range("H91").Value = -228900
$ws.Range("I91").Value = -118300
$ws.Range("J91").Value = -115400
$ws.Range("E102").Value = -7900
$ws.Range("F102").Value = -42400
